$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Chargingdata")

# The first three data rows (5:7) are exact duplicates of rows 2:4 (same
# station rows repeated). Remove the duplicated rows, which shifts every
# row below them up by three and renumbers row references / shared
# formulas accordingly.
$ws.Rows("5:7").Delete()

# Leave the selection where the author's session ended up.
$ws.Range("H12").Select()
